# RTM: updated defect list
# Mark the "In progress" defects in rows 14, 19, 28, 30 and 31 (column F,
# "Status") as "fixed", matching the look already used for the other
# "fixed" rows in the sheet (yellow fill, e.g. row 8/11/12/16...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use an existing "fixed" cell as the formatting template so the new
# cells line up with the rest of the sheet (same fill/border/alignment).
$template = $ws.Range("F8")

$rowsToFix = @(14, 19, 28, 30, 31)
foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "fixed"
    $cell.Interior.Color = $template.Interior.Color
}

# Update the view: scroll/select further down the defect list.
$ws.Activate()
$ws.Range("F31").Select()
